$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 737 (all rows 737..781 shift down to 738..782)
$ws.Rows(737).Insert()

# Populate the newly inserted row 737 with the new record
$ws.Range("A737").Value2 = 6
$ws.Range("B737").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C737").Value2 = "Metropolitana"
$ws.Range("D737").Value2 = 44783
$ws.Range("E737").Value2 = 13
$ws.Range("F737").Value2 = 100112003
$ws.Range("G737").Value2 = "Ajo"
$ws.Range("H737").Value2 = "Chino"
$ws.Range("I737").Value2 = "Primera"
$ws.Range("J737").Value2 = 630
$ws.Range("K737").Value2 = 23000
$ws.Range("L737").Value2 = 24000
$ws.Range("M737").Value2 = 23444
$ws.Range("N737").Value2 = "$/malla 10 kilos"
$ws.Range("O737").Value2 = "China"
$ws.Range("P737").Value2 = 2344
$ws.Range("Q737").Value2 = 10
$ws.Range("R737").Value2 = "Hortaliza"
